# "Update countries & provincias Spain"
#
# The sheet "Pais" (A3:H216) is a COVID-19 snapshot table:
#   A = Pais, B = Casos totales, C = Nuevos casos, D = Casos activos,
#   E = Recuperados, F = Casos criticos, G = Muertes hoy, H = Muertes
# and it is kept sorted by column B (Casos totales) descending.
#
# This refresh updates the per-country counters for a handful of
# countries (Alemania, Austria, Suecia, Dinamarca, Croacia, Uzbekistan,
# Sierra Leona) with newer totals. Because several of those totals moved
# past a neighbouring country's total, the table below also re-sorts the
# affected rows so it stays ordered by "Casos totales".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row([int]$row, [string]$pais, [double]$casosTotales, [double]$nuevosCasos,
                  [double]$casosActivos, [double]$recuperados, [double]$casosCriticos,
                  [double]$muertesHoy, [double]$muertes) {
    $ws.Cells.Item($row, 1).Value = $pais
    $ws.Cells.Item($row, 2).Value = $casosTotales
    $ws.Cells.Item($row, 3).Value = $nuevosCasos
    $ws.Cells.Item($row, 4).Value = $casosActivos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $casosCriticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# --- Straight data refreshes (no change in row order) -----------------
Set-Row 8  "Alemania"  138273 575 81800 52372 4288 49 4101
Set-Row 20 "Austria"   14540  64  9704  4426  227  0  410
Set-Row 32 "Dinamarca" 7073   194 3389  3348  93   15 336
Set-Row 62 "Croacia"   1814   23  600   1178  30   1  36
Set-Row 71 "Uzbekistan" 1390  41  140   1246  8    0  4

# --- Suecia's new total (13216) overtakes Israel (12855): swap rows ---
Set-Row 23 "Suecia" 13216 676 550  11266 482 67 1400
Set-Row 24 "Israel" 12855 97  2967 9740  182 6  148

# --- Sierra Leona jumps from 15 to 26 cases and moves up the table,---
# --- ahead of Zimbabue; every country previously between "Zimbabue"  --
# --- and "Sierra Leona" shifts down one row as a result. --------------
Set-Row 173 "Sierra Leona"        26 11 0  26 0 0 0
Set-Row 174 "Zimbabue"            24 1  2  19 0 0 3
Set-Row 175 "Antigua y Barbuda"   23 0  3  17 1 0 3
Set-Row 176 "Laos"                19 0  2  17 0 0 0
Set-Row 177 "Angola"              19 0  5  12 0 0 2
Set-Row 178 "Timor Oriental"      18 0  1  17 0 0 0
Set-Row 179 "Belice"              18 0  0  16 1 0 2
Set-Row 180 "Nueva Caledonia"     18 0  14 4  1 0 0
Set-Row 181 "Islas Virgenes de los Estados Unidos" 17 0 0 17 0 0 0
Set-Row 182 "Fiyi"                17 0  0  17 0 0 0
Set-Row 183 "Nepal"               16 0  2  14 0 0 0
Set-Row 184 "Malaui"              16 0  0  14 1 0 2
Set-Row 185 "Namibia"             16 0  4  12 0 0 0
Set-Row 186 "Dominica"            16 0  8  8  0 0 0
Set-Row 187 "Suazilandia"         16 0  8  7  0 0 1

# --- Ties re-ordered: San Cristobal y Nieves now precedes Granada -----
Set-Row 190 "San Cristobal y Nieves" 14 0 0 14 0 0 0
Set-Row 191 "Granada"                14 0 0 14 2 0 0

# --- Ties re-ordered: Montserrat now precedes Islas Malvinas ----------
Set-Row 195 "Montserrat"      11 0 1 10 1 0 0
Set-Row 197 "Islas Malvinas"  11 0 1 10 0 0 0

Write-Output "Updated countries & provincias Spain"
